$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "ktorích" -> "ktorých" (the "í" becomes "ý"), which Word records as a
#    run split: "kto" | "rý" | "ch"  (the original run held "ktorí"+"ch" as
#    two runs; editing mid-word splits the touched run further).
# ---------------------------------------------------------------------------
$findRng = $d.Content
$findRng.Find.Execute("ktorích", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$wordStart = $findRng.Start

# temporary bookmarks at the internal run boundaries stop the engine from
# re-coalescing the freshly split runs back together
$b1 = $d.Range($wordStart + 3, $wordStart + 3)
$d.Bookmarks.Add("tmpSplit1", $b1) | Out-Null
$b2 = $d.Range($wordStart + 5, $wordStart + 5)
$d.Bookmarks.Add("tmpSplit2", $b2) | Out-Null

$riRng = $d.Range($wordStart + 3, $wordStart + 5)
$riRng.Text = "rý"

$d.Bookmarks("tmpSplit1").Delete()
$d.Bookmarks("tmpSplit2").Delete()

# ---------------------------------------------------------------------------
# 2) "miezd" stays "miezd", but the cursor's last-edit position (tracked by
#    the hidden "_GoBack" bookmark) moves into the middle of the word,
#    splitting the run into "mi" | "ezd". The bookmark used to sit at the
#    very end of the document - move it here instead.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$miezdRng = $d.Content
$miezdRng.Find.Execute("miezd", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$miezdStart = $miezdRng.Start

# plant the (permanent) "_GoBack" bookmark at the split point *before*
# touching the text - this is what keeps the "mi" run's original rsid
# attribute untouched while "ezd" becomes a brand-new run
$goBackPoint = $d.Range($miezdStart + 2, $miezdStart + 2)
$d.Bookmarks.Add("_GoBack", $goBackPoint) | Out-Null

$ezdRng = $d.Range($miezdStart + 2, $miezdStart + 5)
$ezdRng.Text = "ezX"
$ezXRng = $d.Range($miezdStart + 2, $miezdStart + 5)
$ezXRng.Text = "ezd"
